$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.254443526268005
$ws.Range("B1").Value = 2.231007814407349
$ws.Range("C1").Value = 6.155633926391602
$ws.Range("D1").Value = 1.393369317054749
$ws.Range("E1").Value = 1.34547233581543
